$d = $word.ActiveDocument

# Change 1: fix typo "HLA1-C" -> "HLA-C"
$d.Content.Find.Execute("HLA1-C", $true, $false, $false, $false, $false, $true, 1, $false, "HLA-C", 2)

# Change 2: "However, for disjoint clusters," -> "However, for overlap clusters,"
$d.Content.Find.Execute("However, for disjoint clusters,", $true, $false, $false, $false, $false, $true, 1, $false, "However, for overlap clusters,", 2)
